# Auto-generated edit script
# Applies numeric cell updates to match the target diff across 8 sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 17
$ws.Range("I6").Value = 15.5
$ws.Range("K6").Value = 46.5
$ws.Range("M6").Value = 65.5
$ws.Range("H18").Value = 1269
$ws.Range("I18").Value = 1363.5
$ws.Range("J18").Value = 702
$ws.Range("K18").Value = 1363.5
$ws.Range("L18").Value = 702
$ws.Range("M18").Value = -1079.5
$ws.Range("N18").Value = -1270
$ws.Range("H43").Value = 2000
$ws.Range("J43").Value = 2000
$ws.Range("L43").Value = 2000
$ws.Range("N43").Value = -2138
$ws.Range("H125").Value = 500000220
$ws.Range("J125").Value = 500000000
$ws.Range("L125").Value = 4500000000
$ws.Range("N125").Value = -4500004920

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 898.625
$ws.Range("I2").Value = 669.8570999999999
$ws.Range("K2").Value = 669.8570999999999
$ws.Range("M2").Value = -556.8570999999999
$ws.Range("H116").Value = 898.625
$ws.Range("I116").Value = 669.8570999999999
$ws.Range("K116").Value = 669.8570999999999
$ws.Range("M116").Value = 1624.1429
$ws.Range("H132").Value = 2951.5
$ws.Range("I132").Value = 2951.5
$ws.Range("K132").Value = 8854.5
$ws.Range("M132").Value = -6324.5
$ws.Range("H141").Value = 23500
$ws.Range("J141").Value = 23500
$ws.Range("L141").Value = 23500
$ws.Range("N141").Value = -33860

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 898.625
$ws.Range("I3").Value = 669.8570999999999
$ws.Range("K3").Value = 669.8570999999999
$ws.Range("M3").Value = -555.8570999999999
$ws.Range("H86").Value = 1950
$ws.Range("I86").Value = 1900
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1900
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -777
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 1950
$ws.Range("I89").Value = 1900
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 9500
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -3884
$ws.Range("N89").Value = -21232
$ws.Range("H99").Value = 2095.9092
$ws.Range("I99").Value = 1757
$ws.Range("K99").Value = 1757
$ws.Range("M99").Value = -259

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H60").Value = 4800
$ws.Range("I60").Value = 4800
$ws.Range("K60").Value = 4800
$ws.Range("M60").Value = -4289
$ws.Range("H69").Value = 33972
$ws.Range("I69").Value = 33972
$ws.Range("K69").Value = 33972
$ws.Range("M69").Value = -33223
$ws.Range("H72").Value = 33972
$ws.Range("I72").Value = 33972
$ws.Range("K72").Value = 101916
$ws.Range("M72").Value = -98172
$ws.Range("H122").Value = 3949

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 138300.56
$ws.Range("I2").Value = 183418
$ws.Range("J2").Value = 111230.1
$ws.Range("K2").Value = 1100508
$ws.Range("L2").Value = 667380.6000000001
$ws.Range("M2").Value = -1100395
$ws.Range("N2").Value = -667606.6000000001
$ws.Range("H75").Value = 4001.5715
$ws.Range("I75").Value = 1299
$ws.Range("J75").Value = 6028.5
$ws.Range("K75").Value = 3897
$ws.Range("L75").Value = 18085.5
$ws.Range("M75").Value = -2899
$ws.Range("N75").Value = -20081.5
$ws.Range("H78").Value = 4001.5715
$ws.Range("I78").Value = 1299
$ws.Range("J78").Value = 6028.5
$ws.Range("K78").Value = 11691
$ws.Range("L78").Value = 54256.5
$ws.Range("M78").Value = -6699
$ws.Range("N78").Value = -64240.5
$ws.Range("H114").Value = 1727.8889
$ws.Range("J114").Value = 1422.6923
$ws.Range("L114").Value = 4268.0769
$ws.Range("N114").Value = -10776.0769
$ws.Range("H123").Value = 8996.666999999999
$ws.Range("I123").Value = 7327.3335
$ws.Range("K123").Value = 21982.0005
$ws.Range("M123").Value = -19532.0005
$ws.Range("H130").Value = 1614.4286
$ws.Range("I130").Value = 1606.25
$ws.Range("K130").Value = 4818.75
$ws.Range("M130").Value = 201.25
$ws.Range("H131").Value = 669974.4399999999
$ws.Range("J131").Value = 669974.4399999999
$ws.Range("L131").Value = 2009923.32
$ws.Range("N131").Value = -2020003.32

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4749.5
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 6999
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 20997
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -25897
$ws.Range("H132").Value = 1736.625
$ws.Range("I132").Value = 1636.75
$ws.Range("J132").Value = 1836.5
$ws.Range("K132").Value = 4910.25
$ws.Range("L132").Value = 5509.5
$ws.Range("M132").Value = -2380.25
$ws.Range("N132").Value = -10569.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8024.25
$ws.Range("I7").Value = 4449
$ws.Range("J7").Value = 8535
$ws.Range("K7").Value = 4449
$ws.Range("L7").Value = 8535
$ws.Range("M7").Value = -4337
$ws.Range("N7").Value = -8759
$ws.Range("H46").Value = 2999
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H74").Value = 58999.5
$ws.Range("I74").Value = 58999.5
$ws.Range("K74").Value = 58999.5
$ws.Range("M74").Value = -58001.5
$ws.Range("H76").Value = 17216
$ws.Range("J76").Value = 17216
$ws.Range("L76").Value = 17216
$ws.Range("N76").Value = -17892
$ws.Range("H77").Value = 58999.5
$ws.Range("I77").Value = 58999.5
$ws.Range("K77").Value = 176998.5
$ws.Range("M77").Value = -172006.5
$ws.Range("H79").Value = 17216
$ws.Range("J79").Value = 17216
$ws.Range("L79").Value = 17216
$ws.Range("N79").Value = -19556
$ws.Range("H126").Value = 8024.25
$ws.Range("I126").Value = 4449
$ws.Range("J126").Value = 8535
$ws.Range("K126").Value = 13347
$ws.Range("L126").Value = 25605
$ws.Range("M126").Value = -10877
$ws.Range("N126").Value = -30545

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4339.1113
$ws.Range("I132").Value = 4506.5
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 13519.5
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -10989.5
$ws.Range("N132").Value = -14060
